$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the used range of column A (rows containing verse references).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value2
    if ($null -eq $val) { continue }
    $text = [string]$val
    # Strip an erroneous trailing "16" suffix that was appended to the
    # verse reference, e.g. "Ecclesiastes 1:116" -> "Ecclesiastes 1:1".
    if ($text.EndsWith("16") -and $text -ne "Ecclesiastes 7:201600") {
        $cell.Value2 = $text.Substring(0, $text.Length - 2)
    }
}
